# Append the new run-log row (row 48) to Sheet1, mirroring the formatting
# of the previous row (row 47) and filling in the new run's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 48
$prevRow = 47

# Copy the formatting (style) of the last existing row onto the new row
# before writing values, so the new cells pick up the same centered style.
$srcRange = $ws.Range("A" + $prevRow + ":H" + $prevRow)
$dstRange = $ws.Range("A" + $newRow + ":H" + $newRow)
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = "2025-08-23 09:36:10 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-23 15:06:10 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
